$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 234, pushing the existing rows 234.. down to 236..
# (mirrors a new week of price data being prepended to this product's table)
$ws.Range("A234:A235").EntireRow.Insert()

# Row 234: Betarraga, Primera, new week (2022-03-16 -> serial 44636)
$ws.Range("A234").Value = 7
$ws.Range("B234").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C234").Value = "Ñuble"
$ws.Range("D234").Value = 44636
$ws.Range("E234").Value = 16
$ws.Range("F234").Value = 100114014
$ws.Range("G234").Value = "Betarraga"
$ws.Range("H234").Value = "Sin especificar"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 300
$ws.Range("K234").Value = 700
$ws.Range("L234").Value = 800
$ws.Range("M234").Value = 750
$ws.Range("N234").Value = "$/paquete 5 unidades"
$ws.Range("O234").Value = "Región del Maule"
$ws.Range("P234").Value = 150
$ws.Range("Q234").Value = 5
$ws.Range("R234").Value = "Hortaliza"

# Row 235: Betarraga, Segunda, same new week
$ws.Range("A235").Value = 7
$ws.Range("B235").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C235").Value = "Ñuble"
$ws.Range("D235").Value = 44636
$ws.Range("E235").Value = 16
$ws.Range("F235").Value = 100114014
$ws.Range("G235").Value = "Betarraga"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Segunda"
$ws.Range("J235").Value = 100
$ws.Range("K235").Value = 600
$ws.Range("L235").Value = 600
$ws.Range("M235").Value = 600
$ws.Range("N235").Value = "$/paquete 5 unidades"
$ws.Range("O235").Value = "Región del Maule"
$ws.Range("P235").Value = 120
$ws.Range("Q235").Value = 5
$ws.Range("R235").Value = "Hortaliza"
